$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 17.04324350258874
$ws.Cells.Item(2, 3).Value = 5.59824321883092
$ws.Cells.Item(2, 4).Value = 11.64155665218443
$ws.Cells.Item(2, 5).Value = 11.26414711091088
$ws.Cells.Item(2, 6).Value = 57.93622936338926
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 11).Value = 13.41342992578834
$ws.Cells.Item(2, 12).Value = 10.2068438843021
$ws.Cells.Item(2, 13).Value = 16.55792841314695

# Row 3
$ws.Cells.Item(3, 2).Value = 17.02031697725706
$ws.Cells.Item(3, 3).Value = 5.506379353199279
$ws.Cells.Item(3, 4).Value = 11.5026526335238
$ws.Cells.Item(3, 5).Value = 11.25634870994391
$ws.Cells.Item(3, 6).Value = 56.83751660945542
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 11).Value = 13.4276402657267
$ws.Cells.Item(3, 12).Value = 10.2198928432413
$ws.Cells.Item(3, 13).Value = 16.59012411845326

# Row 4
$ws.Cells.Item(4, 2).Value = 17.01247937329834
$ws.Cells.Item(4, 3).Value = 5.447486261030889
$ws.Cells.Item(4, 4).Value = 11.41534165416855
$ws.Cells.Item(4, 5).Value = 11.25220294936779
$ws.Cells.Item(4, 6).Value = 56.15357939601257
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 11).Value = 13.44154486887152
$ws.Cells.Item(4, 12).Value = 10.22934092039325
$ws.Cells.Item(4, 13).Value = 16.61349014883754

# Row 5
$ws.Cells.Item(5, 2).Value = 17.01085662783807
$ws.Cells.Item(5, 3).Value = 5.422862620463987
$ws.Cells.Item(5, 4).Value = 11.3792659501928
$ws.Cells.Item(5, 5).Value = 11.25067545713132
$ws.Cells.Item(5, 6).Value = 55.87278340861366
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 11).Value = 13.44850991437787
$ws.Cells.Item(5, 12).Value = 10.23355223626528
$ws.Cells.Item(5, 13).Value = 16.62391557373401

# Row 6
$ws.Cells.Item(6, 2).Value = 17.01068209526773
$ws.Cells.Item(6, 3).Value = 5.418736228742604
$ws.Cells.Item(6, 4).Value = 11.37324604819221
$ws.Cells.Item(6, 5).Value = 11.25043160931821
$ws.Cells.Item(6, 6).Value = 55.82603915182026
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 11).Value = 13.44974477585297
$ws.Cells.Item(6, 12).Value = 10.23427333545721
$ws.Cells.Item(6, 13).Value = 16.62570124558379

# Row 7
$ws.Cells.Item(7, 2).Value = 17.01245112535939
$ws.Cells.Item(7, 3).Value = 5.447156703327107
$ws.Cells.Item(7, 4).Value = 11.41485711454424
$ws.Cells.Item(7, 5).Value = 11.2521816928931
$ws.Cells.Item(7, 6).Value = 56.14980058278266
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 11).Value = 13.44163354862225
$ws.Cells.Item(7, 12).Value = 10.22939625335341
$ws.Cells.Item(7, 13).Value = 16.61362709284712

# Row 8
$ws.Cells.Item(8, 2).Value = 17.03404530072985
$ws.Cells.Item(8, 3).Value = 5.567085822274444
$ws.Cells.Item(8, 4).Value = 11.59409021939322
$ws.Cells.Item(8, 5).Value = 11.26132510850138
$ws.Cells.Item(8, 6).Value = 57.55948079948081
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 11).Value = 13.41725279355024
$ws.Cells.Item(8, 12).Value = 10.21104518034357
$ws.Cells.Item(8, 13).Value = 16.56828217153416

# Row 9
$ws.Cells.Item(9, 2).Value = 17.12572081233704
$ws.Cells.Item(9, 3).Value = 5.782472406115036
$ws.Cells.Item(9, 4).Value = 11.92895759623427
$ws.Cells.Item(9, 5).Value = 11.28433671937148
$ws.Cells.Item(9, 6).Value = 60.23896799485446
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 11).Value = 13.41065471049141
$ws.Cells.Item(9, 12).Value = 10.18644900253947
$ws.Cells.Item(9, 13).Value = 16.50795044169869

# Row 10
$ws.Cells.Item(10, 2).Value = 17.2227929691092
$ws.Cells.Item(10, 3).Value = 5.928603314562613
$ws.Cells.Item(10, 4).Value = 12.16414639955409
$ws.Cells.Item(10, 5).Value = 11.30431936471266
$ws.Cells.Item(10, 6).Value = 62.14191684388648
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 11).Value = 13.43102447109717
$ws.Cells.Item(10, 12).Value = 10.17531611461012
$ws.Cells.Item(10, 13).Value = 16.48110324242439

# Row 11
$ws.Cells.Item(11, 2).Value = 17.27329414357646
$ws.Cells.Item(11, 3).Value = 5.992446045993822
$ws.Cells.Item(11, 4).Value = 12.26864198823393
$ws.Cells.Item(11, 5).Value = 11.31407264080171
$ws.Cells.Item(11, 6).Value = 62.99062873151681
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 11).Value = 13.44576468010342
$ws.Cells.Item(11, 12).Value = 10.17175581775988
$ws.Cells.Item(11, 13).Value = 16.47269150197586

# Row 12
$ws.Cells.Item(12, 2).Value = 17.29331749727387
$ws.Cells.Item(12, 3).Value = 6.016242570598912
$ws.Cells.Item(12, 4).Value = 12.30784141913698
$ws.Cells.Item(12, 5).Value = 11.31786073399553
$ws.Cells.Item(12, 6).Value = 63.30936420109349
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 11).Value = 13.45213155604451
$ws.Cells.Item(12, 12).Value = 10.17062362114812
$ws.Cells.Item(12, 13).Value = 16.47005291274653

# Row 13
$ws.Cells.Item(13, 2).Value = 17.28896531240368
$ws.Cells.Item(13, 3).Value = 6.011134443665338
$ws.Cells.Item(13, 4).Value = 12.29941580214371
$ws.Cells.Item(13, 5).Value = 11.31704069693428
$ws.Cells.Item(13, 6).Value = 63.24084012529841
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 11).Value = 13.45072546345438
$ws.Cells.Item(13, 12).Value = 10.17085785784877
$ws.Cells.Item(13, 13).Value = 16.47059686324113

# Row 14
$ws.Cells.Item(14, 2).Value = 17.27492352342281
$ws.Cells.Item(14, 3).Value = 5.994411406238063
$ws.Cells.Item(14, 4).Value = 12.27187444890683
$ws.Cells.Item(14, 5).Value = 11.31438239442544
$ws.Cells.Item(14, 6).Value = 63.01690565868768
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 11).Value = 13.44627276844008
$ws.Cells.Item(14, 12).Value = 10.17165834388433
$ws.Cells.Item(14, 13).Value = 16.47246346587103

# Row 15
$ws.Cells.Item(15, 2).Value = 17.26643928752909
$ws.Cells.Item(15, 3).Value = 5.984118644509532
$ws.Cells.Item(15, 4).Value = 12.25495593059492
$ws.Cells.Item(15, 5).Value = 11.31276642645494
$ws.Cells.Item(15, 6).Value = 62.87938759925649
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 11).Value = 13.44364752440356
$ws.Cells.Item(15, 12).Value = 10.17217678642295
$ws.Cells.Item(15, 13).Value = 16.47367801780897

# Row 16
$ws.Cells.Item(16, 2).Value = 17.21961920555283
$ws.Cells.Item(16, 3).Value = 5.924377994781041
$ws.Cells.Item(16, 4).Value = 12.15726614155825
$ws.Cells.Item(16, 5).Value = 11.30369525667132
$ws.Cells.Item(16, 6).Value = 62.08609170894217
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 11).Value = 13.43017116185267
$ws.Cells.Item(16, 12).Value = 10.17557903236911
$ws.Cells.Item(16, 13).Value = 16.48172947707366

# Row 17
$ws.Cells.Item(17, 2).Value = 17.19251279670343
$ws.Cells.Item(17, 3).Value = 5.887054111484361
$ws.Cells.Item(17, 4).Value = 12.09668923011395
$ws.Cells.Item(17, 5).Value = 11.2982997951807
$ws.Cells.Item(17, 6).Value = 61.59493157174791
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 11).Value = 13.42330473408444
$ws.Cells.Item(17, 12).Value = 10.17805125801337
$ws.Cells.Item(17, 13).Value = 16.4876425413478

# Row 18
$ws.Cells.Item(18, 2).Value = 17.17751975015468
$ws.Cells.Item(18, 3).Value = 5.865338606280474
$ws.Cells.Item(18, 4).Value = 12.06161305155874
$ws.Cells.Item(18, 5).Value = 11.29525891880168
$ws.Cells.Item(18, 6).Value = 61.31084604886988
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 11).Value = 13.41987077039077
$ws.Cells.Item(18, 12).Value = 10.17961481094863
$ws.Cells.Item(18, 13).Value = 16.49140137106505

# Row 19
$ws.Cells.Item(19, 2).Value = 17.17254639251608
$ws.Cells.Item(19, 3).Value = 5.857943575222732
$ws.Cells.Item(19, 4).Value = 12.04969706201321
$ws.Cells.Item(19, 5).Value = 11.29424007413719
$ws.Cells.Item(19, 6).Value = 61.21439444614592
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 11).Value = 13.41879665836613
$ws.Cells.Item(19, 12).Value = 10.18016852960079
$ws.Cells.Item(19, 13).Value = 16.49273548884596

# Row 20
$ws.Cells.Item(20, 2).Value = 17.19533653039916
$ws.Cells.Item(20, 3).Value = 5.891052954026034
$ws.Cells.Item(20, 4).Value = 12.10316204801346
$ws.Cells.Item(20, 5).Value = 11.29886769215235
$ws.Cells.Item(20, 6).Value = 61.64738186772283
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 11).Value = 13.42398234417829
$ws.Cells.Item(20, 12).Value = 10.17777343278359
$ws.Cells.Item(20, 13).Value = 16.48697605377444

# Row 21
$ws.Cells.Item(21, 2).Value = 17.2790236278719
$ws.Cells.Item(21, 3).Value = 5.999333667492725
$ws.Cells.Item(21, 4).Value = 12.27997417072181
$ws.Cells.Item(21, 5).Value = 11.31516063617305
$ws.Cells.Item(21, 6).Value = 63.0827543081442
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 11).Value = 13.44755934803393
$ws.Cells.Item(21, 12).Value = 10.17141736176449
$ws.Cells.Item(21, 13).Value = 16.47190036066097

# Row 22
$ws.Cells.Item(22, 2).Value = 17.33895495605296
$ws.Cells.Item(22, 3).Value = 6.067890129010989
$ws.Cells.Item(22, 4).Value = 12.39336452217614
$ws.Cells.Item(22, 5).Value = 11.32636091869166
$ws.Cells.Item(22, 6).Value = 64.00529653413358
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 11).Value = 13.46754258186486
$ws.Cells.Item(22, 12).Value = 10.16852220789094
$ws.Cells.Item(22, 13).Value = 16.46523432773525

# Row 23
$ws.Cells.Item(23, 2).Value = 17.30649381559622
$ws.Cells.Item(23, 3).Value = 6.031502774657141
$ws.Cells.Item(23, 4).Value = 12.3330481516595
$ws.Cells.Item(23, 5).Value = 11.32033282975741
$ws.Cells.Item(23, 6).Value = 63.51440988542013
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 11).Value = 13.45645957347572
$ws.Cells.Item(23, 12).Value = 10.16995232248023
$ws.Cells.Item(23, 13).Value = 16.46850053422873

# Row 24
$ws.Cells.Item(24, 2).Value = 17.19405807921835
$ws.Cells.Item(24, 3).Value = 5.889245879401487
$ws.Cells.Item(24, 4).Value = 12.10023646396074
$ws.Cells.Item(24, 5).Value = 11.29861075575513
$ws.Cells.Item(24, 6).Value = 61.62367440177412
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 11).Value = 13.42367439656816
$ws.Cells.Item(24, 12).Value = 10.17789859450033
$ws.Cells.Item(24, 13).Value = 16.48727625353632

# Row 25
$ws.Cells.Item(25, 2).Value = 17.09566864287698
$ws.Cells.Item(25, 3).Value = 5.72632924776329
$ws.Cells.Item(25, 4).Value = 11.8402199783667
$ws.Cells.Item(25, 5).Value = 11.27756955653331
$ws.Cells.Item(25, 6).Value = 59.52466198294143
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 11).Value = 13.40801379871023
$ws.Cells.Item(25, 12).Value = 10.19188392218123
$ws.Cells.Item(25, 13).Value = 16.52120511107043
